$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "Dimension before:" $ws.Range("A1").Worksheet.UsedRange.Address()
$ws.Rows("665:666").Insert()
Write-Host "Inserted"
